# Apply updated DAMSLTag (column I) and DialogAct (column J) values
# following re-run of SGNN dialog act annotation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$changes = @(
    @{Row=15; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=51; I='sv'; J='Statement-opinion'}
    @{Row=56; I='%'; J='Uninterpretable'}
    @{Row=87; I='ba'; J='Appreciation'}
    @{Row=93; I='ba'; J='Appreciation'}
    @{Row=114; I='aa'; J='Agree/Accept'}
    @{Row=115; I='aa'; J='Agree/Accept'}
    @{Row=124; I='sd'; J='Statement-non-opinion'}
    @{Row=125; I='sd'; J='Statement-non-opinion'}
    @{Row=126; I='sd'; J='Statement-non-opinion'}
    @{Row=128; I='qy'; J='Yes-No-Question'}
    @{Row=129; I='sv'; J='Statement-opinion'}
    @{Row=133; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=149; I='aa'; J='Agree/Accept'}
    @{Row=166; I='aa'; J='Agree/Accept'}
    @{Row=167; I='sd'; J='Statement-non-opinion'}
    @{Row=170; I='aa'; J='Agree/Accept'}
    @{Row=174; I='ba'; J='Appreciation'}
    @{Row=181; I='sv'; J='Statement-opinion'}
    @{Row=184; I='aa'; J='Agree/Accept'}
    @{Row=186; I='sv'; J='Statement-opinion'}
    @{Row=190; I='qy'; J='Yes-No-Question'}
    @{Row=196; I='sd'; J='Statement-non-opinion'}
    @{Row=199; I='sd'; J='Statement-non-opinion'}
    @{Row=220; I='sd'; J='Statement-non-opinion'}
    @{Row=229; I='aa'; J='Agree/Accept'}
    @{Row=239; I='sd'; J='Statement-non-opinion'}
    @{Row=253; I='%'; J='Uninterpretable'}
    @{Row=271; I='sd'; J='Statement-non-opinion'}
    @{Row=279; I='%'; J='Uninterpretable'}
    @{Row=280; I='aa'; J='Agree/Accept'}
    @{Row=282; I='aa'; J='Agree/Accept'}
    @{Row=296; I='aa'; J='Agree/Accept'}
    @{Row=306; I='sd'; J='Statement-non-opinion'}
    @{Row=307; I='ba'; J='Appreciation'}
    @{Row=313; I='aa'; J='Agree/Accept'}
    @{Row=316; I='sd'; J='Statement-non-opinion'}
    @{Row=319; I='aa'; J='Agree/Accept'}
    @{Row=322; I='sd'; J='Statement-non-opinion'}
    @{Row=336; I='sd'; J='Statement-non-opinion'}
    @{Row=346; I='sv'; J='Statement-opinion'}
    @{Row=347; I='ba'; J='Appreciation'}
    @{Row=364; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=366; I='sv'; J='Statement-opinion'}
    @{Row=373; I='sv'; J='Statement-opinion'}
    @{Row=377; I='sd'; J='Statement-non-opinion'}
    @{Row=380; I='qy'; J='Yes-No-Question'}
    @{Row=387; I='sd'; J='Statement-non-opinion'}
    @{Row=392; I='sd'; J='Statement-non-opinion'}
    @{Row=417; I='sv'; J='Statement-opinion'}
    @{Row=433; I='%'; J='Uninterpretable'}
    @{Row=434; I='aa'; J='Agree/Accept'}
    @{Row=435; I='%'; J='Uninterpretable'}
    @{Row=436; I='sv'; J='Statement-opinion'}
    @{Row=439; I='sv'; J='Statement-opinion'}
    @{Row=446; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=448; I='sd'; J='Statement-non-opinion'}
    @{Row=452; I='aa'; J='Agree/Accept'}
    @{Row=456; I='b'; J='Acknowledge (Backchannel)'}
)

foreach ($change in $changes) {
    $ws.Cells.Item($change.Row, 9).Value = $change.I
    $ws.Cells.Item($change.Row, 10).Value = $change.J
}
